$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.395.57"
$ws.Range("E2").Value = "  +0.00%  "

$ws.Range("D3").Value = "1.846.71"
$ws.Range("E3").Value = "  -0.21%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.9973"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.23%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "240.39"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.01%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.6276"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.16%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.9993"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.15%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.07493"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.73%  "

$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("E10").Value = "  -1.18%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07736"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.09%  "

$ws.Range("D12").Value = "1.847.13"
$ws.Range("E12").Value = "  -2.18%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "4.999"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.74%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.6808"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.27%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.00001055"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.20%  "

$ws.Range("E16").Value = "  -1.19%  "

$ws.Range("D17").Value = "2.108.77"
$ws.Range("E17").Value = "  -3.57%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "6.183"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.17%  "

$ws.Range("D19").Value = "29.442.40"
$ws.Range("E19").Value = "  +0.11%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "229.82"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.76%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "12.33"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.24%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.9989"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.12%  "

$ws.Range("E23").Value = "  -0.16%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.9993"
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "159.00"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.11%  "

$ws.Range("E26").Value = "  -0.79%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "8.421"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.14%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "17.55"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.94%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.06515"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +16.20%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.417"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +1.35%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.478"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.32%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "4.100"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.32%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.103"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.83%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.831"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.28%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.142"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.94%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.6997"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.21%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.578"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.10%  "

$ws.Range("D38").Value = "1.270.93"
$ws.Range("E38").Value = "  +2.87%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.831"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +4.22%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.01836"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.56%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.9084"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.51%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.9985"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.24%  "

$ws.Range("D44").Value = "2.011.31"
$ws.Range("E44").Value = "  -18.27%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "101.43"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.10%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "66.37"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.45%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.749"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +4.11%  "

$ws.Range("E48").Value = "  +1.22%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "7.079"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.01%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.1174"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.54%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "9.076"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.94%  "
